$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 315 (appended): constant columns copied from the pattern used by all data rows
$ws.Cells.Item(315, 1).Value = 3
$ws.Cells.Item(315, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(315, 3).Value = 'Coquimbo'
$ws.Cells.Item(315, 5).Value = 5
$ws.Cells.Item(315, 6).Value = 100112039
$ws.Cells.Item(315, 7).Value = 'Ciboulette'
$ws.Cells.Item(315, 8).Value = 'Sin especificar'
$ws.Cells.Item(315, 9).Value = 'Primera'
$ws.Cells.Item(315, 14).Value = '$/docena de atados'
$ws.Cells.Item(315, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(315, 17).Value = 3
$ws.Cells.Item(315, 18).Value = 'Hortaliza'

# Match the date-formatted style used by every other row's Fecha (column D) cell.
$ws.Cells.Item(315, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Weekly Fecha/Volumen/Precio series (cols D,J,K,L,M,P) shift down by one row for rows
# 16..315 (each row inherits what used to be directly above it), and row 15 starts a new
# entry at the head of the series.
$series = @(
    ,(15, 44699, 160, 1500, 1500, 1500, 500)
    ,(16, 44179, 48, 2000, 2000, 2000, 667)
    ,(17, 44662, 160, 1500, 1500, 1500, 500)
    ,(18, 44176, 80, 1500, 1500, 1500, 500)
    ,(19, 44210, 120, 1500, 1500, 1500, 500)
    ,(20, 44516, 150, 1500, 1500, 1500, 500)
    ,(21, 44435, 810, 1500, 1500, 1500, 500)
    ,(22, 44391, 160, 1500, 1500, 1500, 500)
    ,(23, 44231, 120, 1500, 1500, 1500, 500)
    ,(24, 44364, 160, 1500, 1500, 1500, 500)
    ,(25, 44320, 160, 1500, 1500, 1500, 500)
    ,(26, 44336, 160, 1500, 1500, 1500, 500)
    ,(27, 44425, 160, 1500, 1500, 1500, 500)
    ,(28, 44496, 150, 1500, 1500, 1500, 500)
    ,(29, 44522, 160, 1500, 1500, 1500, 500)
    ,(30, 44630, 160, 1500, 1500, 1500, 500)
    ,(31, 44299, 130, 1500, 1500, 1500, 500)
    ,(32, 44175, 120, 1500, 1500, 1500, 500)
    ,(33, 44509, 160, 1500, 1500, 1500, 500)
    ,(34, 44265, 120, 1500, 1500, 1500, 500)
    ,(35, 44512, 160, 1500, 1500, 1500, 500)
    ,(36, 44638, 180, 1500, 1500, 1500, 500)
    ,(37, 44397, 160, 1500, 1500, 1500, 500)
    ,(38, 44433, 180, 1500, 1500, 1500, 500)
    ,(39, 44592, 160, 1500, 1500, 1500, 500)
    ,(40, 44414, 160, 1500, 1500, 1500, 500)
    ,(41, 44168, 160, 1500, 1500, 1500, 500)
    ,(42, 44473, 160, 1500, 1500, 1500, 500)
    ,(43, 44460, 160, 1500, 1500, 1500, 500)
    ,(44, 44355, 180, 1500, 1500, 1500, 500)
    ,(45, 44657, 120, 1500, 1500, 1500, 500)
    ,(46, 44613, 120, 1500, 1500, 1500, 500)
    ,(47, 44623, 160, 1500, 1500, 1500, 500)
    ,(48, 44482, 160, 1500, 1500, 1500, 500)
    ,(49, 44305, 180, 1500, 1500, 1500, 500)
    ,(50, 44216, 80, 1500, 1500, 1500, 500)
    ,(51, 44195, 180, 1500, 1500, 1500, 500)
    ,(52, 44690, 160, 1500, 1500, 1500, 500)
    ,(53, 44271, 180, 1500, 1500, 1500, 500)
    ,(54, 44386, 160, 1500, 1500, 1500, 500)
    ,(55, 44426, 160, 1500, 1500, 1500, 500)
    ,(56, 44313, 130, 1500, 1500, 1500, 500)
    ,(57, 44186, 180, 1500, 1500, 1500, 500)
    ,(58, 44543, 160, 1500, 1500, 1500, 500)
    ,(59, 44407, 160, 1500, 1500, 1500, 500)
    ,(60, 44489, 160, 1500, 1500, 1500, 500)
    ,(61, 44284, 180, 1500, 1500, 1500, 500)
    ,(62, 44235, 160, 1500, 1500, 1500, 500)
    ,(63, 44319, 190, 1500, 1500, 1500, 500)
    ,(64, 44253, 120, 1500, 1500, 1500, 500)
    ,(65, 44578, 250, 1500, 1500, 1500, 500)
    ,(66, 44504, 160, 1500, 1500, 1500, 500)
    ,(67, 44194, 80, 1500, 1500, 1500, 500)
    ,(68, 44188, 180, 1500, 1500, 1500, 500)
    ,(69, 44204, 180, 1500, 1500, 1500, 500)
    ,(70, 44356, 160, 1500, 1500, 1500, 500)
    ,(71, 44285, 160, 1500, 1500, 1500, 500)
    ,(72, 44214, 110, 1500, 1500, 1500, 500)
    ,(73, 44392, 160, 1500, 1500, 1500, 500)
    ,(74, 44540, 180, 1500, 1500, 1500, 500)
    ,(75, 44438, 160, 1500, 1500, 1500, 500)
    ,(76, 44428, 160, 1500, 1500, 1500, 500)
    ,(77, 44559, 172, 1500, 2000, 1747, 582)
    ,(78, 44466, 160, 1500, 1500, 1500, 500)
    ,(79, 44557, 80, 1500, 1500, 1500, 500)
    ,(80, 44560, 180, 1500, 1500, 1500, 500)
    ,(81, 44208, 160, 1500, 1500, 1500, 500)
    ,(82, 44258, 230, 1500, 1500, 1500, 500)
    ,(83, 44218, 130, 1500, 1500, 1500, 500)
    ,(84, 44406, 160, 1500, 1500, 1500, 500)
    ,(85, 44614, 230, 1500, 1500, 1500, 500)
    ,(86, 44302, 130, 1500, 1500, 1500, 500)
    ,(87, 44193, 120, 1800, 1800, 1800, 600)
    ,(88, 44203, 120, 1500, 1500, 1500, 500)
    ,(89, 44524, 160, 1500, 1500, 1500, 500)
    ,(90, 44455, 160, 1500, 1500, 1500, 500)
    ,(91, 44484, 160, 1500, 1500, 1500, 500)
    ,(92, 44274, 120, 1500, 1500, 1500, 500)
    ,(93, 44461, 160, 1500, 1500, 1500, 500)
    ,(94, 44624, 160, 1500, 1500, 1500, 500)
    ,(95, 44573, 160, 1500, 1500, 1500, 500)
    ,(96, 44608, 120, 1500, 1500, 1500, 500)
    ,(97, 44161, 180, 1500, 1500, 1500, 500)
    ,(98, 44617, 160, 1500, 1500, 1500, 500)
    ,(99, 44312, 160, 1500, 1500, 1500, 500)
    ,(100, 44547, 160, 1500, 1500, 1500, 500)
    ,(101, 44615, 160, 1500, 1500, 1500, 500)
    ,(102, 44211, 120, 1500, 1500, 1500, 500)
    ,(103, 44264, 120, 1500, 1500, 1500, 500)
    ,(104, 44434, 140, 1500, 1500, 1500, 500)
    ,(105, 44379, 160, 1500, 1500, 1500, 500)
    ,(106, 44229, 160, 1500, 1500, 1500, 500)
    ,(107, 44550, 160, 1500, 1500, 1500, 500)
    ,(108, 44419, 130, 1500, 1500, 1500, 500)
    ,(109, 44278, 130, 1500, 1500, 1500, 500)
    ,(110, 44663, 120, 1500, 1500, 1500, 500)
    ,(111, 44546, 180, 1500, 1500, 1500, 500)
    ,(112, 44659, 120, 1500, 1500, 1500, 500)
    ,(113, 44293, 160, 1500, 1500, 1500, 500)
    ,(114, 44596, 160, 1500, 1500, 1500, 500)
    ,(115, 44649, 160, 1500, 1500, 1500, 500)
    ,(116, 44529, 160, 1500, 1500, 1500, 500)
    ,(117, 44467, 160, 1500, 1500, 1500, 500)
    ,(118, 44385, 180, 1500, 1500, 1500, 500)
    ,(119, 44532, 160, 1500, 1500, 1500, 500)
    ,(120, 44306, 160, 1500, 1500, 1500, 500)
    ,(121, 44636, 160, 1500, 1500, 1500, 500)
    ,(122, 44172, 110, 1500, 1500, 1500, 500)
    ,(123, 44580, 160, 1500, 1500, 1500, 500)
    ,(124, 44420, 160, 1500, 1500, 1500, 500)
    ,(125, 44665, 120, 1500, 1500, 1500, 500)
    ,(126, 44487, 160, 1500, 1500, 1500, 500)
    ,(127, 44272, 160, 1500, 1500, 1500, 500)
    ,(128, 44645, 120, 1500, 1500, 1500, 500)
    ,(129, 44236, 120, 1500, 1500, 1500, 500)
    ,(130, 44286, 160, 1500, 1500, 1500, 500)
    ,(131, 44308, 160, 1500, 1500, 1500, 500)
    ,(132, 44595, 230, 1500, 1500, 1500, 500)
    ,(133, 44494, 190, 1500, 1500, 1500, 500)
    ,(134, 44421, 180, 1500, 1500, 1500, 500)
    ,(135, 44403, 180, 1500, 1500, 1500, 500)
    ,(136, 44432, 150, 1500, 1500, 1500, 500)
    ,(137, 44553, 150, 1500, 1500, 1500, 500)
    ,(138, 44295, 120, 1500, 1500, 1500, 500)
    ,(139, 44687, 120, 1500, 1500, 1500, 500)
    ,(140, 44181, 90, 1500, 1500, 1500, 500)
    ,(141, 44307, 130, 1500, 1500, 1500, 500)
    ,(142, 44651, 120, 1500, 1500, 1500, 500)
    ,(143, 44642, 160, 1500, 1500, 1500, 500)
    ,(144, 44476, 160, 1500, 1500, 1500, 500)
    ,(145, 44526, 160, 1500, 1500, 1500, 500)
    ,(146, 44561, 180, 1500, 1500, 1500, 500)
    ,(147, 44329, 160, 1500, 1500, 1500, 500)
    ,(148, 44545, 180, 1500, 1500, 1500, 500)
    ,(149, 44518, 160, 1500, 1500, 1500, 500)
    ,(150, 44348, 160, 1500, 1500, 1500, 500)
    ,(151, 44446, 180, 1500, 1500, 1500, 500)
    ,(152, 44350, 160, 1500, 1500, 1500, 500)
    ,(153, 44530, 120, 1500, 1500, 1500, 500)
    ,(154, 44398, 160, 1500, 1500, 1500, 500)
    ,(155, 44200, 120, 1500, 1500, 1500, 500)
    ,(156, 44431, 180, 1500, 1500, 1500, 500)
    ,(157, 44567, 180, 1500, 1500, 1500, 500)
    ,(158, 44322, 130, 1500, 1500, 1500, 500)
    ,(159, 44370, 180, 1500, 1500, 1500, 500)
    ,(160, 44327, 190, 1500, 1500, 1500, 500)
    ,(161, 44246, 160, 1500, 1500, 1500, 500)
    ,(162, 44452, 190, 1500, 1500, 1500, 500)
    ,(163, 44259, 120, 1500, 1500, 1500, 500)
    ,(164, 44300, 160, 1500, 1500, 1500, 500)
    ,(165, 44383, 180, 1500, 1500, 1500, 500)
    ,(166, 44321, 130, 1500, 1500, 1500, 500)
    ,(167, 44362, 180, 1500, 1500, 1500, 500)
    ,(168, 44266, 120, 1500, 1500, 1500, 500)
    ,(169, 44517, 160, 1500, 1500, 1500, 500)
    ,(170, 44607, 120, 1500, 1500, 1500, 500)
    ,(171, 44237, 130, 1500, 1500, 1500, 500)
    ,(172, 44468, 180, 1500, 1500, 1500, 500)
    ,(173, 44539, 160, 1500, 1500, 1500, 500)
    ,(174, 44162, 160, 1500, 1500, 1500, 500)
    ,(175, 44270, 120, 1500, 1500, 1500, 500)
    ,(176, 44643, 120, 1500, 1500, 1500, 500)
    ,(177, 44589, 150, 1500, 1500, 1500, 500)
    ,(178, 44483, 180, 1500, 1500, 1500, 500)
    ,(179, 44196, 180, 1500, 1500, 1500, 500)
    ,(180, 44396, 160, 1500, 1500, 1500, 500)
    ,(181, 44249, 160, 1500, 1500, 1500, 500)
    ,(182, 44497, 160, 1500, 1500, 1500, 500)
    ,(183, 44252, 160, 1500, 1500, 1500, 500)
    ,(184, 44620, 120, 1500, 1500, 1500, 500)
    ,(185, 44202, 120, 1500, 1500, 1500, 500)
    ,(186, 44453, 130, 1500, 1500, 1500, 500)
    ,(187, 44588, 180, 1500, 1500, 1500, 500)
    ,(188, 44298, 160, 1500, 1500, 1500, 500)
    ,(189, 44581, 130, 1500, 1500, 1500, 500)
    ,(190, 44503, 160, 1500, 1500, 1500, 500)
    ,(191, 44616, 160, 1500, 1500, 1500, 500)
    ,(192, 44628, 130, 1500, 1500, 1500, 500)
    ,(193, 44244, 110, 1500, 1500, 1500, 500)
    ,(194, 44585, 160, 1500, 1500, 1500, 500)
    ,(195, 44335, 160, 1500, 1500, 1500, 500)
    ,(196, 44341, 160, 1500, 1500, 1500, 500)
    ,(197, 44277, 160, 1500, 1500, 1500, 500)
    ,(198, 44441, 190, 1500, 1500, 1500, 500)
    ,(199, 44554, 120, 1500, 1500, 1500, 500)
    ,(200, 44491, 160, 1500, 1500, 1500, 500)
    ,(201, 44544, 160, 1500, 1500, 1500, 500)
    ,(202, 44685, 120, 1500, 1500, 1500, 500)
    ,(203, 44328, 160, 1500, 1500, 1500, 500)
    ,(204, 44201, 120, 1500, 1500, 1500, 500)
    ,(205, 44418, 150, 1500, 1500, 1500, 500)
    ,(206, 44609, 120, 1500, 1500, 1500, 500)
    ,(207, 44410, 120, 1500, 1500, 1500, 500)
    ,(208, 44384, 160, 1500, 1500, 1500, 500)
    ,(209, 44330, 160, 1500, 1500, 1500, 500)
    ,(210, 44399, 120, 1500, 1500, 1500, 500)
    ,(211, 44232, 120, 1500, 1500, 1500, 500)
    ,(212, 44677, 120, 1500, 1500, 1500, 500)
    ,(213, 44238, 130, 1500, 1500, 1500, 500)
    ,(214, 44250, 160, 1500, 1500, 1500, 500)
    ,(215, 44334, 190, 1500, 1500, 1500, 500)
    ,(216, 44280, 120, 1500, 1500, 1500, 500)
    ,(217, 44622, 120, 1500, 1500, 1500, 500)
    ,(218, 44454, 160, 1500, 1500, 1500, 500)
    ,(219, 44635, 160, 1500, 1500, 1500, 500)
    ,(220, 44637, 160, 1500, 1500, 1500, 500)
    ,(221, 44582, 180, 1500, 1500, 1500, 500)
    ,(222, 44333, 120, 1500, 1500, 1500, 500)
    ,(223, 44301, 130, 1500, 1500, 1500, 500)
    ,(224, 44698, 120, 1500, 1500, 1500, 500)
    ,(225, 44257, 120, 1500, 1500, 1500, 500)
    ,(226, 44495, 160, 1500, 1500, 1500, 500)
    ,(227, 44498, 160, 1500, 1500, 1500, 500)
    ,(228, 44314, 160, 1500, 1500, 1500, 500)
    ,(229, 44209, 160, 1500, 1500, 1500, 500)
    ,(230, 44189, 180, 1500, 1500, 1500, 500)
    ,(231, 44217, 120, 1500, 1500, 1500, 500)
    ,(232, 44349, 160, 1500, 1500, 1500, 500)
    ,(233, 44215, 130, 1500, 1500, 1500, 500)
    ,(234, 44627, 45, 2000, 2000, 2000, 667)
    ,(235, 44405, 160, 1500, 1500, 1500, 500)
    ,(236, 44650, 110, 1500, 1500, 1500, 500)
    ,(237, 44475, 160, 1500, 1500, 1500, 500)
    ,(238, 44239, 120, 1500, 1500, 1500, 500)
    ,(239, 44358, 160, 1500, 1500, 1500, 500)
    ,(240, 44273, 160, 1500, 1500, 1500, 500)
    ,(241, 44658, 180, 1500, 1500, 1500, 500)
    ,(242, 44400, 160, 1500, 1500, 1500, 500)
    ,(243, 44382, 160, 1500, 1500, 1500, 500)
    ,(244, 44167, 150, 1500, 1500, 1500, 500)
    ,(245, 44551, 160, 1500, 1500, 1500, 500)
    ,(246, 44281, 160, 1500, 1500, 1500, 500)
    ,(247, 44412, 160, 1500, 1500, 1500, 500)
    ,(248, 44357, 160, 1500, 1500, 1500, 500)
    ,(249, 44363, 130, 1500, 1500, 1500, 500)
    ,(250, 44558, 160, 1500, 1500, 1500, 500)
    ,(251, 44694, 160, 1500, 1500, 1500, 500)
    ,(252, 44459, 160, 1500, 1500, 1500, 500)
    ,(253, 44552, 180, 1500, 1500, 1500, 500)
    ,(254, 44488, 150, 1500, 1500, 1500, 500)
    ,(255, 44316, 160, 1500, 1500, 1500, 500)
    ,(256, 44568, 160, 1500, 1500, 1500, 500)
    ,(257, 44656, 160, 1500, 1500, 1500, 500)
    ,(258, 44586, 160, 1500, 1500, 1500, 500)
    ,(259, 44469, 160, 1500, 1500, 1500, 500)
    ,(260, 44463, 160, 1500, 1500, 1500, 500)
    ,(261, 44579, 160, 1500, 1500, 1500, 500)
    ,(262, 44342, 260, 1500, 1500, 1500, 500)
    ,(263, 44160, 230, 1500, 1500, 1500, 500)
    ,(264, 44251, 80, 1500, 1500, 1500, 500)
    ,(265, 44279, 160, 1500, 1500, 1500, 500)
    ,(266, 44610, 160, 1500, 1500, 1500, 500)
    ,(267, 44519, 160, 1500, 1500, 1500, 500)
    ,(268, 44505, 120, 1500, 1500, 1500, 500)
    ,(269, 44372, 160, 1500, 1500, 1500, 500)
    ,(270, 44267, 160, 1500, 1500, 1500, 500)
    ,(271, 44669, 85, 2000, 2000, 2000, 667)
    ,(272, 44477, 160, 1500, 1500, 1500, 500)
    ,(273, 44671, 150, 1500, 1500, 1500, 500)
    ,(274, 44474, 160, 1500, 1500, 1500, 500)
    ,(275, 44571, 190, 1500, 1500, 1500, 500)
    ,(276, 44631, 160, 1500, 1500, 1500, 500)
    ,(277, 44490, 160, 1500, 1500, 1500, 500)
    ,(278, 44679, 180, 1500, 1500, 1500, 500)
    ,(279, 44369, 180, 1500, 1500, 1500, 500)
    ,(280, 44365, 180, 1500, 1500, 1500, 500)
    ,(281, 44603, 160, 1500, 1500, 1500, 500)
    ,(282, 44575, 120, 1500, 1500, 1500, 500)
    ,(283, 44427, 160, 1500, 1500, 1500, 500)
    ,(284, 44587, 120, 1500, 1500, 1500, 500)
    ,(285, 44565, 180, 1500, 1500, 1500, 500)
    ,(286, 44447, 160, 1500, 1500, 1500, 500)
    ,(287, 44445, 180, 1500, 1500, 1500, 500)
    ,(288, 44533, 160, 1500, 1500, 1500, 500)
    ,(289, 44523, 160, 1500, 1500, 1500, 500)
    ,(290, 44601, 160, 1500, 1500, 1500, 500)
    ,(291, 44343, 180, 1500, 1500, 1500, 500)
    ,(292, 44159, 120, 1500, 1500, 1500, 500)
    ,(293, 44629, 130, 1500, 1500, 1500, 500)
    ,(294, 44606, 160, 1500, 1500, 1500, 500)
    ,(295, 44594, 130, 1500, 1500, 1500, 500)
    ,(296, 44377, 160, 1500, 1500, 1500, 500)
    ,(297, 44417, 160, 1500, 1500, 1500, 500)
    ,(298, 44566, 130, 1500, 1500, 1500, 500)
    ,(299, 44344, 160, 1500, 1500, 1500, 500)
    ,(300, 44351, 160, 1500, 1500, 1500, 500)
    ,(301, 44508, 160, 1500, 1500, 1500, 500)
    ,(302, 44600, 160, 1500, 1500, 1500, 500)
    ,(303, 44323, 160, 1500, 1500, 1500, 500)
    ,(304, 44515, 160, 1500, 1500, 1500, 500)
    ,(305, 44602, 130, 1500, 1500, 1500, 500)
    ,(306, 44326, 120, 1500, 1500, 1500, 500)
    ,(307, 44692, 160, 1500, 1500, 1500, 500)
    ,(308, 44165, 68, 2000, 2000, 2000, 667)
    ,(309, 44655, 120, 1500, 1500, 1500, 500)
    ,(310, 44315, 130, 1500, 1500, 1500, 500)
    ,(311, 44448, 160, 1500, 1500, 1500, 500)
    ,(312, 44263, 180, 1500, 1500, 1500, 500)
    ,(313, 44648, 180, 1500, 1500, 1500, 500)
    ,(314, 44376, 160, 1500, 1500, 1500, 500)
    ,(315, 44442, 180, 1500, 1500, 1500, 500)
)

foreach ($entry in $series) {
    $row = $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
    $ws.Cells.Item($row, 11).Value = $entry[3]
    $ws.Cells.Item($row, 12).Value = $entry[4]
    $ws.Cells.Item($row, 13).Value = $entry[5]
    $ws.Cells.Item($row, 16).Value = $entry[6]
}
